# Update automatico via Actualizar 06-20-2020 05-34-57
# Appends the next day's COVID-19 patient-condition record to the
# "Condicion_Pacientes" table on sheet "Hoja1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data row (date 2020-06-19, serial 44001) appended right below the
# existing table range (A1:F98) so the ListObject auto-expands to A1:F99.
# Carry the previous row's formatting down, the way Excel does when a
# table grows (new row inherits the style of the row above it).
$ws.Range("A98:F98").Copy()
$ws.Range("A99:F99").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A99").Value = 44001
$ws.Range("B99").Value = 1056
$ws.Range("C99").Value = 519
$ws.Range("D99").Value = 603
$ws.Range("E99").Value = 318
$ws.Range("F99").Value = 46

# The worksheet range grew by one row; expand the Excel Table
# ("Condicion_Pacientes") to include it, same as Excel does automatically
# when you type into the row right below a table.
$table = $ws.ListObjects.Item("Condicion_Pacientes")
$table.Resize($ws.Range("A1:F99"))

# Match the saved selection/view state from the source workbook.
$ws.Range("F99").Select()
